$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "21_FL_TTCA" column (C) for "1_Antihypertensive"
$ws.Columns("C").Insert()
$ws.Range("C1").Value = "1_Antihypertensive"
$ws.Range("C2").Value = "non-active"
$ws.Range("C3").Value = "non-active"

# Insert a new column before the existing "23_FL_AMAP_alternative" column (now E) for "22_FL_umami"
$ws.Columns("E").Insert()
$ws.Range("E1").Value = "22_FL_umami"
$ws.Range("E2").Value = "non-active"
$ws.Range("E3").Value = "non-active"

# Append two new columns after "24_FL_AMAP_main" (G): "25_FL_AMP" and "26_FL_MRSA".
# Copy the header formatting from the existing "24_FL_AMAP_main" header cell (G1)
# so the new header cells (H1, I1) keep the same bold/border/centered style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("I1").PasteSpecial(-4122)

$ws.Range("H1").Value = "25_FL_AMP"
$ws.Range("H2").Value = "non-active"
$ws.Range("H3").Value = "non-active"

$ws.Range("I1").Value = "26_FL_MRSA"
$ws.Range("I2").Value = "active"
$ws.Range("I3").Value = "active"
